$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: classical-best-embeddings vs. classical-best-tfidf
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.057
$ws.Range("D2").Value = 0.022
$ws.Range("E2").Value = 0.02
$ws.Range("F2").Value = 0.018
$ws.Range("H2").Value = 0.038
$ws.Range("I2").Value = 0.029
$ws.Range("J2").Value = 0.033

# Row 3: BERT-base vs. classical-best-tfidf (name unchanged)
$ws.Range("C3").Value = 0.074
$ws.Range("D3").Value = 0.079
$ws.Range("E3").Value = 0.076
$ws.Range("F3").Value = 0.059
$ws.Range("G3").Value = 0.112
$ws.Range("H3").Value = 0.115
$ws.Range("J3").Value = 0.08599999999999999

# Row 4: BERT-base vs. classical-best-embeddings
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.017
$ws.Range("D4").Value = 0.057
$ws.Range("E4").Value = 0.056
$ws.Range("F4").Value = 0.041
$ws.Range("G4").Value = 0.07000000000000001
$ws.Range("H4").Value = 0.077
$ws.Range("I4").Value = 0.043
$ws.Range("J4").Value = 0.053

# Row 5: BERT-base-nli vs. classical-best-tfidf (name unchanged)
$ws.Range("B5").Value = 0.528
$ws.Range("C5").Value = 0.114
$ws.Range("D5").Value = 0.08599999999999999
$ws.Range("E5").Value = 0.077
$ws.Range("F5").Value = 0.055
$ws.Range("G5").Value = 0.093
$ws.Range("H5").Value = 0.093
$ws.Range("I5").Value = 0.083
$ws.Range("J5").Value = 0.08599999999999999

# Row 6: BERT-base-nli vs. classical-best-embeddings
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.528
$ws.Range("C6").Value = 0.057
$ws.Range("D6").Value = 0.064
$ws.Range("E6").Value = 0.057
$ws.Range("F6").Value = 0.037
$ws.Range("G6").Value = 0.051
$ws.Range("H6").Value = 0.055
$ws.Range("I6").Value = 0.054
$ws.Range("J6").Value = 0.054

# Row 7: BERT-base-nli vs. BERT-base (name unchanged)
$ws.Range("B7").Value = 0.528
$ws.Range("C7").Value = 0.04
$ws.Range("D7").Value = 0.007
$ws.Range("E7").Value = 0.001
$ws.Range("F7").Value = -0.004
$ws.Range("G7").Value = -0.019
$ws.Range("H7").Value = -0.022
$ws.Range("I7").Value = 0.011
$ws.Range("J7").Value = 0.001
